$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cache price/volume columns as Text so numeric-looking strings
# (e.g. "1.005", "26.028.14") are preserved exactly as typed.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.028.14'
$ws.Range("D3").Value = '1.668.03'
$ws.Range("E3").Value = '  -1.43%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '217.04'
$ws.Range("E5").Value = '  -1.26%  '
$ws.Range("D6").Value = '0.5104'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '0.2659'
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.06403'
$ws.Range("E9").Value = '  +1.16%  '
$ws.Range("D10").Value = '21.78'
$ws.Range("E10").Value = '  -1.52%  '
$ws.Range("D11").Value = '0.07452'
$ws.Range("E11").Value = '  +1.53%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '4.514'
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.665.77'
$ws.Range("E13").Value = '  -1.66%  '
$ws.Range("D14").Value = '0.5834'
$ws.Range("E14").Value = '  +0.73%  '
$ws.Range("D15").Value = '0.000008566'
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").Value = '64.30'
$ws.Range("E16").Value = '  -1.72%  '
$ws.Range("D17").Value = '26.089.84'
$ws.Range("E17").Value = '  -2.05%  '
$ws.Range("D18").Value = '4.940'
$ws.Range("E18").Value = '  -0.97%  '
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("D20").Value = '10.78'
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("D21").Value = '191.65'
$ws.Range("E21").Value = '  +2.57%  '
$ws.Range("D22").Value = '6.206'
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").Value = '144.87'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").Value = '7.618'
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("D26").Value = '0.1201'
$ws.Range("E26").Value = '  +2.64%  '
$ws.Range("D27").Value = '15.66'
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("D28").Value = '0.06527'
$ws.Range("E28").Value = '  +13.56%  '
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("D30").Value = '1.317'
$ws.Range("E30").Value = '  -1.82%  '
$ws.Range("D31").Value = '3.541'
$ws.Range("E31").Value = '  +0.52%  '
$ws.Range("D32").Value = '3.519'
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("D33").Value = '1.649'
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("D34").Value = '1.020'
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").Value = '0.6118'
$ws.Range("E35").Value = '  +1.95%  '
$ws.Range("D36").Value = '2.371'
$ws.Range("E36").Value = '  +0.26%  '
$ws.Range("D37").Value = '2.682'
$ws.Range("E37").Value = '  -0.30%  '
$ws.Range("D38").Value = '6.255'
$ws.Range("E38").Value = '  +7.22%  '
$ws.Range("D39").Value = '0.01602'
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("D40").Value = '1.092.13'
$ws.Range("E40").Value = '  -0.50%  '
$ws.Range("D41").Value = '0.8652'
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("D43").Value = '101.10'
$ws.Range("E43").Value = '  +1.55%  '
$ws.Range("D44").Value = '1.815.97'
$ws.Range("E44").Value = '  -1.89%  '
$ws.Range("D46").Value = '56.44'
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").Value = '1.006'
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("D48").Value = '8.073'
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("D50").Value = '6.095'
$ws.Range("E50").Value = '  +5.22%  '
$ws.Range("D51").Value = '0.4286'
$ws.Range("E51").Value = '  -0.90%  '

# Restore default (unstyled) cell style now that text is committed,
# matching the workbook's original unstyled data cells.
$priceRange.Style = "Normal"
